$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the first 29 data rows (rows 2-30), shifting everything else up.
$ws.Range("A2:A30").EntireRow.Delete()

# Leave the active cell where Excel naturally lands after this kind of edit.
$ws.Range("D9").Select()
